# Up to date as of 2/17
# Reschedules the "Meeting # N" / "Meeting with Advisor # N" rows of the
# Gantt chart (rows 21-36), marks the now-past meetings (through 2/17) as
# "Completed", drops "Meeting with Advisor # 7" and appends a new
# "Meeting # 14" entry at the bottom of the schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: Online Café Talk -> now Completed (dates unchanged) ---
$ws.Range("F21").Value = "Completed"

# --- Row 22: Meeting with Advisor # 2, moved to 2/13, now Completed ---
$ws.Range("B22").Value = "Meeting with Advisor # 2"
$ws.Range("C22").Value = 43144
$ws.Range("D22").Value = 43144
$ws.Range("F22").Value = "Completed"

# --- Row 23: Meeting # 7, moved to 2/17, now Completed ---
$ws.Range("B23").Value = "Meeting # 7"
$ws.Range("C23").Value = 43148
$ws.Range("D23").Value = 43148
$ws.Range("F23").Value = "Completed"

# Rows 21-23 pick up the "Completed" status look (no fill), matching the
# style already used further up the sheet for completed tasks.
$ws.Range("F8").Copy()
$ws.Range("F21:F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 24: Meeting # 8, moved to 2/22 ---
$ws.Range("C24").Value = 43153
$ws.Range("D24").Value = 43153

# --- Row 25: Meeting # 9, moved to 2/25 ---
$ws.Range("B25").Value = "Meeting # 9"
$ws.Range("C25").Value = 43156
$ws.Range("D25").Value = 43156

# --- Row 26: Meeting with Advisor # 3, moved to 2/27 ---
$ws.Range("B26").Value = "Meeting with Advisor # 3"
$ws.Range("C26").Value = 43158
$ws.Range("D26").Value = 43158

# --- Row 27: Meeting # 10, moved to 3/3 ---
$ws.Range("B27").Value = "Meeting # 10"
$ws.Range("C27").Value = 43162
$ws.Range("D27").Value = 43162

# --- Row 28: Meeting with Advisor # 4, moved to 3/6 ---
$ws.Range("B28").Value = "Meeting with Advisor # 4"
$ws.Range("C28").Value = 43165
$ws.Range("D28").Value = 43165

# --- Row 29: Meeting # 11, moved to 3/8 ---
$ws.Range("B29").Value = "Meeting # 11"
$ws.Range("C29").Value = 43167
$ws.Range("D29").Value = 43167

# --- Row 30: Meeting with Advisor # 5, moved to 3/10 (dates already matched) ---
$ws.Range("B30").Value = "Meeting with Advisor # 5"

# --- Row 31: Meeting # 12, moved to 3/11 (single day, was a 6-day span) ---
$ws.Range("B31").Value = "Meeting # 12"
$ws.Range("C31").Value = 43170
$ws.Range("D31").Value = 43170

# --- Row 32: Meeting # 13, moved to 3/13 (dates already matched) ---
$ws.Range("B32").Value = "Meeting # 13"

# --- Row 33: First Paper Draft Review with Advisor, moved to 3/13 ---
$ws.Range("B33").Value = "First Paper Draft Review with Advisor"
$ws.Range("C33").Value = 43172
$ws.Range("D33").Value = 43172

# --- Row 35: was "Meeting with Advisor # 7" -> now "Meeting with Advisor # 6" ---
$ws.Range("B35").Value = "Meeting with Advisor # 6"

# --- Row 36: was "Meeting # 13" -> now "Meeting # 14" (new meeting added) ---
$ws.Range("B36").Value = "Meeting # 14"
